$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# (wordRow, wordCol, newText) -- wordRow/wordCol are 1-indexed table cell coordinates
$replacements = @(
    @(1, 1, "51÷4=12, 3"),
    @(1, 2, "81÷4=20, 1"),
    @(1, 3, "13÷9=1, 4"),
    @(1, 4, "91÷3=30, 1"),
    @(1, 5, "31÷9=3, 4"),
    @(5, 1, "40÷9=4, 4"),
    @(5, 2, "26÷4=6, 2"),
    @(5, 3, "16÷9=1, 7"),
    @(5, 4, "11÷3=3, 2"),
    @(5, 5, "28÷8=3, 4"),
    @(9, 1, "78÷3=26, 0"),
    @(9, 3, "48÷4=12, 0"),
    @(9, 4, "66÷5=13, 1"),
    @(9, 5, "35÷6=5, 5"),
    @(13, 1, "67÷5=13, 2"),
    @(13, 2, "61÷3=20, 1"),
    @(13, 3, "92÷2=46, 0"),
    @(13, 4, "81÷4=20, 1"),
    @(13, 5, "34÷7=4, 6"),
    @(17, 1, "96÷5=19, 1"),
    @(17, 2, "58÷7=8, 2"),
    @(17, 3, "74÷4=18, 2"),
    @(17, 4, "67÷3=22, 1"),
    @(17, 5, "45÷9=5, 0")
)

foreach ($item in $replacements) {
    $r = $item[0]
    $c = $item[1]
    $new = $item[2]
    $tbl.Cell($r, $c).Range.Text = $new
}
